# Commit: && - 19387 от 16.06.2025 https://2eurostore.ru/ - 20133 от 20.07.2025 https://2eurostore.ru/
# Adds a new 2025 Belgian commemorative 2-euro coin entry (row 39 on the "2€" sheet):
#   "Circuit de Spa-Francorchamps", Staff of Mercury mint mark, new map of Europe,
#   Aster flower mint director mark, mintage 155000, NL-side "Can exchange" flag set to 1.
# Also restores the last-used selections on both sheets.

$wb = $excel.ActiveWorkbook

$wsCoins = $wb.Worksheets.Item("2€")
$wsLinks = $wb.Worksheets.Item("Links")

# --- Row 39: new commemorative coin data ---
$wsCoins.Range("B39").Value = "Circuit de Spa-Francorchamps"
$wsCoins.Range("D39").Value = "Obv: With mint symbol - `nStaff of Mercury"
$wsCoins.Range("E39").Value = "Rev: new map of Europe"
$wsCoins.Range("F39").Value = "Obv: Mint director symbol - Aster flower"
$wsCoins.Range("I39").Value = 155000
$wsCoins.Range("K39").Value = 1

# Writing the multi-line value into D39 auto-fits the row taller; restore the
# original fixed row height (15) to match the source formatting.
$wsCoins.Rows.Item(39).RowHeight = 15

# --- Restore view selections (active cell) on both sheets ---
$wsCoins.Activate()
$wsCoins.Range("G26").Select()

$wsLinks.Activate()
$wsLinks.Range("B7").Select()

# Leave the coins sheet as the active/selected one, matching the source file
$wsCoins.Activate()
